$d = $word.ActiveDocument

# =====================================================================
# Paragraph 1 - Title: "Abstract" -> "ABSTRACT", sz 36 -> 28, add szCs=28,
# add single underline.
# =====================================================================
$titleRng = $d.Paragraphs(1).Range
$titleRng.Font.Size = 14        # half-points *2 => sz="28"
$titleRng.Font.SizeBi = 14      # => szCs="28"
$titleRng.Font.Underline = 1    # => <w:u w:val="single"/>
$titleRng.Text = "ABSTRACT"

Write-Output "Paragraph 1 done"

# =====================================================================
# Paragraphs 2 & 3 - set base font size 24 (12pt) on the whole paragraph
# (this also stamps the pPr/rPr paragraph-mark run properties).
# =====================================================================
$d.Paragraphs(2).Range.Font.Size = 12
$d.Paragraphs(3).Range.Font.Size = 12

Write-Output "Base sizes done"

# =====================================================================
# Paragraph 2 - insert new opening sentence before "Cyberbullying refers..."
# (merge into the existing run's text so formatting - incl. w:cs - is
# inherited, then split the leading chunk back out with its own Font.Size).
# =====================================================================
$openRng = $d.Content
$openRng.Find.Execute("Cyberbullying refers to the use of technology") | Out-Null
$openInsertStart = $openRng.Start
$openNewText = "Social networking sites are primarily used for communicating and connecting with other people. However, some people use these technologies to harm others emotionally. "
$openRng.Text = $openNewText + "Cyberbullying refers to the use of technology"
$openFormatRng = $d.Range($openInsertStart, $openInsertStart + $openNewText.Length)
$openFormatRng.Font.Size = 11
$openFormatRng.Font.Size = 12

Write-Output "Opening sentence inserted"

# =====================================================================
# Paragraph 2 - replace everything from "Recently it has become..." to
# the end of the paragraph (just before the _GoBack bookmark) with the
# new closing content.
# =====================================================================
$tailFind = $d.Content
$tailFind.Find.Execute("Recently it has become") | Out-Null
$tailStart = $tailFind.Start
$p2Rng = $d.Paragraphs(2).Range
$tailEnd = $p2Rng.End - 1   # stop right before the paragraph mark / bookmark
$tailRange = $d.Range($tailStart, $tailEnd)

$tailNewText = "Truly, cyberbullying has become rampant in many countries. In order to address this issue, the team proposed a method to detect cyberbullying entries on social media through the use of Natural Language Processing (NLP). The data that was used by the researchers was collected from Youtube, Twitter, and Facebook. Each statement that was extracted from these social networking sites were ranked according to their harmfulness level (wherein 0, indicates no indication of cyberbullying event, 1 indicates mild cyberbullying occurrences and 2 indicates severe cyberbullying occurrences). Furthermore, it was annotated into six cyberbullying categories: bad description, intelligence, physical appearance, race and culture, sexuality, and social rejection. Among the 450 words that were extracted from 600 statements, bad description was most prevalent with a frequency of 27%. "

$tailRange.Text = $tailNewText
$tailFormatRange = $d.Range($tailStart, $tailStart + $tailNewText.Length)
$tailFormatRange.Font.Size = 11
$tailFormatRange.Font.Size = 12

Write-Output "Tail replaced"
